$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.725.76"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.889.83"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.67"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4759"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2939"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06538"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.05"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07741"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7412"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.83"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.886.83"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.254"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "276.26"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.725.16"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.22"
$ws.Range("E18").Value = "  -2.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007566"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.135.43"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.334"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9992"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.247"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.251"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.00"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.85"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.346"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09727"
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.305"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.196"
$ws.Range("E33").Value = "  +2.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04890"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.127"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7010"
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.722"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01918"
$ws.Range("E38").Value = "  +2.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.798"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.354"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.95"
$ws.Range("E41").Value = "  +6.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.034"
$ws.Range("E42").Value = "  +3.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4261"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8437"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.45"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.422"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.073"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.73"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "919.63"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05778"
$ws.Range("E51").Value = "  +2.34%  "
